$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final roster table (header in row 1 is unchanged).
# Columns: A = Oyuncu Adı (Player), B = Pozisyon (Position), C = Takım (Team)
$data = @(
    @("D'Angelo Russell", "PG", "Brooklyn Nets"),
    @("Norman Powell", "SG,SF", "LA Clippers"),
    @("Jalen Brunson", "PG", "New York Knicks"),
    @("Devin Booker", "PG,SG", "Phoenix Suns"),
    @("Bilal Coulibaly", "SG,SF", "Washington Wizards"),
    @("Kawhi Leonard", "SG,SF,PF", "LA Clippers"),
    @("Josh Okogie", "SG,SF", "Charlotte Hornets"),
    @("Christian Braun", "SG,SF", "Denver Nuggets"),
    @("Jalen Williams", "SG,SF,PF,C", "Oklahoma City Thunder"),
    @("Myles Turner", "C", "Indiana Pacers"),
    @("LeBron James", "SF,PF", "Los Angeles Lakers"),
    @("Trae Young", "PG", "Atlanta Hawks"),
    @("Scoot Henderson", "PG", "Portland Trail Blazers"),
    @("Walker Kessler", "C", "Utah Jazz"),
    @("Desmond Bane", "SG,SF", "Memphis Grizzlies"),
    @("Immanuel Quickley", "PG,SG", "Toronto Raptors"),
    @("Brandon Ingram", "SG,SF,PF", "New Orleans Pelicans"),
    @("Jimmy Butler", "SF,PF", "Miami Heat")
)

$row = 2
foreach ($rec in $data) {
    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]
    $row = $row + 1
}
